# Add more sample data and improve researcher card layout
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Researchers
$ws2 = $wb.Worksheets.Item(2)   # Projects

# --- Researchers sheet: tweak Diego's description, add two new researchers ---
$ws1.Range("C3").Value = "Supports applied AI prototyping and research web tools focused on hospital operations."

$ws1.Range("A4").Value = "Maya Singh"
$ws1.Range("B4").Value = "Data Scientist"
$ws1.Range("C4").Value = "Builds demand forecasting pipelines and evaluates model performance across hospital units."
$ws1.Range("D4").Value = "https://placehold.co/400x400/png?text=Maya"
$ws1.Hyperlinks.Add($ws1.Range("D4"), "https://placehold.co/400x400/png?text=Maya")
$ws1.Range("D4").Style = "Hyperlink"
$ws1.Range("E4").Value = "maya-singh"

$ws1.Range("A5").Value = "Jordan Lee"
$ws1.Range("B5").Value = "Health Systems Engineer"
$ws1.Range("C5").Value = "Partners with clinical teams to translate analytics into staffing and workflow improvements."
$ws1.Range("D5").Value = "https://placehold.co/400x400/png?text=Jordan"
$ws1.Hyperlinks.Add($ws1.Range("D5"), "https://placehold.co/400x400/png?text=Jordan")
$ws1.Range("D5").Style = "Hyperlink"
$ws1.Range("E5").Value = "jordan-lee"

# --- Projects sheet: update Predict's slide deck text, add two new projects ---
$ws2.Range("G2").Value = "Overview::https://placehold.co/960x540/png?text=Predict+Slide+1::High-level overview|Challenges::https://placehold.co/960x540/png?text=Predict+Slide+2::Key system constraints"

$ws2.Range("A3").Value = "nursing-demand-forecast"
$ws2.Range("B3").Value = "Nursing Demand Forecast"
$ws2.Range("C3").Value = "Forecast nurse demand and optimize staffing plans across critical units."
$ws2.Range("D3").Value = "https://placehold.co/900x600/png?text=Nursing+Demand+Forecast"
$ws2.Hyperlinks.Add($ws2.Range("D3"), "https://placehold.co/900x600/png?text=Nursing+Demand+Forecast")
$ws2.Range("D3").Style = "Hyperlink"
$ws2.Range("E3").Value = "Maya Singh::Data Scientist::maya-singh|Jordan Lee::Health Systems Engineer::jordan-lee"
$ws2.Range("F3").Value = "Integrates historical staffing and census data.|Produces 12-month forecasts by unit.|Flags high-risk gaps for mitigation planning."
$ws2.Range("G3").Value = "Model Inputs::https://placehold.co/960x540/png?text=Demand+Slide+1::Historical and demographic signals|Results::https://placehold.co/960x540/png?text=Demand+Slide+2::Projected utilization curves"

$ws2.Range("A4").Value = "attrition-phenotyping"
$ws2.Range("B4").Value = "Attrition Phenotyping"
$ws2.Range("C4").Value = "Identify distinct attrition profiles and drivers across nursing populations."
$ws2.Range("D4").Value = "https://placehold.co/900x600/png?text=Attrition+Phenotyping"
$ws2.Hyperlinks.Add($ws2.Range("D4"), "https://placehold.co/900x600/png?text=Attrition+Phenotyping")
$ws2.Range("D4").Style = "Hyperlink"
$ws2.Range("E4").Value = "Maya Singh::Data Scientist::maya-singh|Theofilos::Principal Investigator::theofilos"
$ws2.Range("F4").Value = "Clusters separations by role and tenure.|Combines well-being surveys with HR data.|Surfaces modifiable risk factors for intervention."
$ws2.Range("G4").Value = "Cluster View::https://placehold.co/960x540/png?text=Attrition+Slide+1::Phenotype clusters|Interventions::https://placehold.co/960x540/png?text=Attrition+Slide+2::Mitigation strategies"
